# Populate the "Peak Counts" data for both sheets (epidermis + dermis).
# Each sheet already has a header row (Unmanipulated | NoMAP | MAP) in A1:C1.
# We fill rows 2-19 with the counted peak values; some columns run shorter
# than others, leaving trailing cells blank for that row.

$wb = $excel.ActiveWorkbook

$epidermis = $wb.Worksheets.Item("epidermis")
$dermis    = $wb.Worksheets.Item("dermis")

$epidermisData = @(
    @(7, 3, 4),
    @(5, 7, 9),
    @(18, 3, 5),
    @(4, 4, 3),
    @(4, 3, 7),
    @($null, 5, 3),
    @($null, 3, 2),
    @($null, $null, 4),
    @($null, $null, 4),
    @($null, $null, 4),
    @($null, $null, 4),
    @($null, $null, 5),
    @($null, $null, 4),
    @($null, $null, 4),
    @($null, $null, 16),
    @($null, $null, 4),
    @($null, $null, 8),
    @($null, $null, 8)
)

$dermisData = @(
    @(2, 5, 3),
    @(5, 4, 4),
    @(1, 5, 3),
    @(4, 4, 4),
    @(3, 3, 4),
    @($null, 3, 4),
    @($null, 3, 2),
    @($null, $null, 4),
    @($null, $null, 4),
    @($null, $null, 3),
    @($null, $null, 4),
    @($null, $null, 2),
    @($null, $null, 3),
    @($null, $null, 4),
    @($null, $null, 5),
    @($null, $null, 6),
    @($null, $null, 6),
    @($null, $null, 5)
)

$row = 2
foreach ($rowValues in $epidermisData) {
    for ($col = 1; $col -le 3; $col++) {
        $value = $rowValues[$col - 1]
        if ($null -ne $value) {
            $epidermis.Cells.Item($row, $col).Value = $value
        }
    }
    $row++
}

$row = 2
foreach ($rowValues in $dermisData) {
    for ($col = 1; $col -le 3; $col++) {
        $value = $rowValues[$col - 1]
        if ($null -ne $value) {
            $dermis.Cells.Item($row, $col).Value = $value
        }
    }
    $row++
}
